# Update cryptos list values (Price and Volume(1h) columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.107.56"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.15"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.47"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6289"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07613"
$ws.Range("E8").Value = "  +3.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2937"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.69"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.96"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.974"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6667"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001007"
$ws.Range("E15").Value = "  +16.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.04"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.075"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.105.07"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.94"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.41"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.232"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.65"
$ws.Range("E24").Value = "  +0.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.533"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1387"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.500"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.113"
$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.027"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.196"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05272"
$ws.Range("E32").Value = "  -0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.846"
$ws.Range("E33").Value = "  +0.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7366"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.244.42"
$ws.Range("E37").Value = "  -3.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.765"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.377"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8987"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.19"
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.986.97"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000125"
$ws.Range("E45").Value = "  +4.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.53"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5119"
$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4051"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.930"
$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05762"
$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.701"
$ws.Range("E51").Value = "  +0.27%  "
